# Daily attendance processing - 2025-10-09 21:18:11
#
# Normalises the "Recorded By" column (G): for every row whose recorded-by
# list mentions "System" (any case) among its comma-separated entries, the
# order of the entries is reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$col = 7  # column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text.ToLower().Contains("system")) {
        $parts = $text -split ", "
        $count = $parts.Count

        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newText = $reversedParts -join ", "
        $cell.Value = $newText
    }
}
